# Natmi following Dr Hou advice
#
# The underlying NATMI ligand/receptor run was redone: the "M2" sending
# cluster now aggregates 3 ligand-expressing cells (was 1) which changes all
# the derived statistics for the existing M2 -> {ECs, FAPs, M2, sCs} edges in
# rows 2-5, and a brand new "sCs" sending cluster (2 cells) contributes four
# more edges -> {ECs, FAPs, M2, sCs} appended as rows 6-9.
#
# Sending/ligand/receptor columns (A-D) for rows 2-5 are textually unchanged
# (still M2 / Areg / Egfr / <target>), so only the numeric columns E-T need to
# be rewritten there. Rows 6-9 are entirely new.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: M2, Areg, Egfr -> ECs -------------------------------------------------
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.277672333333333
$ws.Range("H2").Value = 3.833017
$ws.Range("I2").Value = 0.7545320937513288
$ws.Range("J2").Value = 0.7545320937513288
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.701354
$ws.Range("N2").Value = 8.104061999999999
$ws.Range("O2").Value = 0.02221077311549548
$ws.Range("P2").Value = 0.02221077311549548
$ws.Range("Q2").Value = 3.451445268339333
$ws.Range("R2").Value = 31.063007415054
$ws.Range("S2").Value = 0.01675874114267053
$ws.Range("T2").Value = 0.01675874114267053

# --- Row 3: M2, Areg, Egfr -> FAPs ------------------------------------------------
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.277672333333333
$ws.Range("H3").Value = 3.833017
$ws.Range("I3").Value = 0.7545320937513288
$ws.Range("J3").Value = 0.7545320937513288
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("N3").Value = 264.449356
$ws.Range("O3").Value = 0.7247753838328104
$ws.Range("P3").Value = 0.7247753838328105
$ws.Range("Q3").Value = 112.6265419096724
$ws.Range("R3").Value = 1013.638877187052
$ws.Range("S3").Value = 0.5468662878627933
$ws.Range("T3").Value = 0.5468662878627935

# --- Row 4: M2, Areg, Egfr -> M2 --------------------------------------------------
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.277672333333333
$ws.Range("H4").Value = 3.833017
$ws.Range("I4").Value = 0.7545320937513288
$ws.Range("J4").Value = 0.7545320937513288
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.24063
$ws.Range("N4").Value = 0.72189
$ws.Range("O4").Value = 0.001978481285600361
$ws.Range("P4").Value = 0.001978481285600361
$ws.Range("Q4").Value = 0.30744629357
$ws.Range("R4").Value = 2.76701664213
$ws.Range("S4").Value = 0.001492827626871861
$ws.Range("T4").Value = 0.001492827626871861

# --- Row 5: M2, Areg, Egfr -> sCs -------------------------------------------------
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.277672333333333
$ws.Range("H5").Value = 3.833017
$ws.Range("I5").Value = 0.7545320937513288
$ws.Range("J5").Value = 0.7545320937513288
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.53182233333333
$ws.Range("N5").Value = 91.595467
$ws.Range("O5").Value = 0.2510353617660938
$ws.Range("P5").Value = 0.2510353617660938
$ws.Range("Q5").Value = 39.00966468154878
$ws.Range("R5").Value = 351.086982133939
$ws.Range("S5").Value = 0.189414237118993
$ws.Range("T5").Value = 0.189414237118993

# --- Row 6 (new): sCs, Areg, Egfr -> ECs ------------------------------------------
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Areg"
$ws.Range("C6").Value = "Egfr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4156583333333333
$ws.Range("H6").Value = 1.246975
$ws.Range("I6").Value = 0.2454679062486713
$ws.Range("J6").Value = 0.2454679062486713
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.701354
$ws.Range("N6").Value = 8.104061999999999
$ws.Range("O6").Value = 0.02221077311549548
$ws.Range("P6").Value = 0.02221077311549548
$ws.Range("Q6").Value = 1.122840301383333
$ws.Range("R6").Value = 10.10556271245
$ws.Range("S6").Value = 0.005452031972824953
$ws.Range("T6").Value = 0.005452031972824953

# --- Row 7 (new): sCs, Areg, Egfr -> FAPs -----------------------------------------
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Areg"
$ws.Range("C7").Value = "Egfr"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4156583333333333
$ws.Range("H7").Value = 1.246975
$ws.Range("I7").Value = 0.2454679062486713
$ws.Range("J7").Value = 0.2454679062486713
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("N7").Value = 264.449356
$ws.Range("O7").Value = 0.7247753838328104
$ws.Range("P7").Value = 0.7247753838328105
$ws.Range("Q7").Value = 36.64019285534444
$ws.Range("R7").Value = 329.7617356980999
$ws.Range("S7").Value = 0.177909095970017
$ws.Range("T7").Value = 0.1779090959700171

# --- Row 8 (new): sCs, Areg, Egfr -> M2 -------------------------------------------
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Areg"
$ws.Range("C8").Value = "Egfr"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.4156583333333333
$ws.Range("H8").Value = 1.246975
$ws.Range("I8").Value = 0.2454679062486713
$ws.Range("J8").Value = 0.2454679062486713
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.24063
$ws.Range("N8").Value = 0.72189
$ws.Range("O8").Value = 0.001978481285600361
$ws.Range("P8").Value = 0.001978481285600361
$ws.Range("Q8").Value = 0.10001986475
$ws.Range("R8").Value = 0.90017878275
$ws.Range("S8").Value = 0.0004856536587285001
$ws.Range("T8").Value = 0.0004856536587285001

# --- Row 9 (new): sCs, Areg, Egfr -> sCs ------------------------------------------
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Areg"
$ws.Range("C9").Value = "Egfr"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.4156583333333333
$ws.Range("H9").Value = 1.246975
$ws.Range("I9").Value = 0.2454679062486713
$ws.Range("J9").Value = 0.2454679062486713
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 30.53182233333333
$ws.Range("N9").Value = 91.595467
$ws.Range("O9").Value = 0.2510353617660938
$ws.Range("P9").Value = 0.2510353617660938
$ws.Range("Q9").Value = 12.69080638470278
$ws.Range("R9").Value = 114.217257462325
$ws.Range("S9").Value = 0.06162112464710079
$ws.Range("T9").Value = 0.06162112464710079
